$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.459.31'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.573.33'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").Value = '''288.06'
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("D7").Value = '''0.3718'
$ws.Range("E7").Value = '  +2.08%  '
$ws.Range("D8").Value = '''47.51'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").Value = '''0.3317'
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("D10").Value = '''1.152'
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("D11").Value = '''0.07528'
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = '''20.77'
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").Value = '''5.940'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").Value = '''6.934'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '1.575.99'
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").Value = '''0.00001119'
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").Value = '''88.19'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '''0.06729'
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = '''6.405'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '''16.52'
$ws.Range("E22").Value = '  +2.75%  '
$ws.Range("D23").Value = '''12.02'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").Value = '22.441.93'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '''2.387'
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("D26").Value = '''2.624'
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("D27").Value = '''151.27'
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("D28").Value = '''19.65'
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").Value = '''4.943'
$ws.Range("E29").Value = '  -1.36%  '
$ws.Range("D30").Value = '''125.26'
$ws.Range("E30").Value = '  +1.79%  '
$ws.Range("D31").Value = '1.750.83'
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("D32").Value = '''1.096'
$ws.Range("E32").Value = '  +3.24%  '
$ws.Range("D33").Value = '''6.089'
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").Value = '''9.855'
$ws.Range("E35").Value = '  +2.60%  '
$ws.Range("D36").Value = '''0.08358'
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").Value = '''0.02455'
$ws.Range("E37").Value = '  +2.57%  '
$ws.Range("D38").Value = '''0.2235'
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("D40").Value = '''0.06383'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '''5.346'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '''11.42'
$ws.Range("E42").Value = '  +2.25%  '
$ws.Range("D43").Value = '''0.6278'
$ws.Range("E43").Value = '  +3.25%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '''1.000'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''13.98'
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("D46").Value = '''0.6094'
$ws.Range("E46").Value = '  +6.01%  '
$ws.Range("D47").Value = '''3.774'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").Value = '''2.050'
$ws.Range("E48").Value = '  +2.08%  '
$ws.Range("D49").Value = '''125.11'
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("D50").Value = '''1.209'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").Value = '''0.07208'
$ws.Range("E51").Value = '  -0.07%  '
